## screenClient.xlsx -- formDef.json cannot have integers in the value list.
## Convert every numeric "value" in the choices sheet's column B to a text
## string (prefixed with "a"), and move the active tab / selection from the
## "survey" sheet to the "choices" sheet (landing on B6), matching the author's
## edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# --- Convert numeric choice "value" cells in column B to text ("a" + N) ---
# (left untouched: rows already stored as text, e.g. B54:B58, B84:B85)
$ws.Range("B2").Value  = "a1"
$ws.Range("B3").Value  = "a0"
$ws.Range("B4").Value  = "a1"
$ws.Range("B5").Value  = "a0"
$ws.Range("B6").Value  = "a888"
$ws.Range("B7").Value  = "a888"
$ws.Range("B8").Value  = "a1"
$ws.Range("B9").Value  = "a2"
$ws.Range("B10").Value = "a3"
$ws.Range("B11").Value = "a4"
$ws.Range("B12").Value = "a5"
$ws.Range("B13").Value = "a6"
$ws.Range("B14").Value = "a7"
$ws.Range("B15").Value = "a888"
$ws.Range("B16").Value = "a1"
$ws.Range("B17").Value = "a2"
$ws.Range("B18").Value = "a3"
$ws.Range("B19").Value = "a4"
$ws.Range("B20").Value = "a5"
$ws.Range("B21").Value = "a6"
$ws.Range("B22").Value = "a7"
$ws.Range("B23").Value = "a8"
$ws.Range("B24").Value = "a9"
$ws.Range("B25").Value = "a10"
$ws.Range("B26").Value = "a888"
$ws.Range("B27").Value = "a0"
$ws.Range("B28").Value = "a1"
$ws.Range("B29").Value = "a2"
$ws.Range("B30").Value = "a3"
$ws.Range("B31").Value = "a4"
$ws.Range("B32").Value = "a5"
$ws.Range("B33").Value = "a6"
$ws.Range("B34").Value = "a7"
$ws.Range("B35").Value = "a8"
$ws.Range("B36").Value = "a9"
$ws.Range("B37").Value = "a888"
$ws.Range("B38").Value = "a1"
$ws.Range("B39").Value = "a2"
$ws.Range("B40").Value = "a3"
$ws.Range("B41").Value = "a888"
$ws.Range("B42").Value = "a0"
$ws.Range("B43").Value = "a1"
$ws.Range("B44").Value = "a888"
$ws.Range("B45").Value = "a999"
$ws.Range("B46").Value = "a1"
$ws.Range("B47").Value = "a0"
$ws.Range("B48").Value = "a888"
$ws.Range("B49").Value = "a999"
$ws.Range("B50").Value = "a1"
$ws.Range("B51").Value = "a2"
$ws.Range("B52").Value = "a3"
$ws.Range("B53").Value = "a4"
$ws.Range("B59").Value = "a0"
$ws.Range("B60").Value = "a1"
$ws.Range("B61").Value = "a2"
$ws.Range("B62").Value = "a3"
$ws.Range("B63").Value = "a4"
$ws.Range("B64").Value = "a5"
$ws.Range("B65").Value = "a888"
$ws.Range("B66").Value = "a9999"
$ws.Range("B67").Value = "a0"
$ws.Range("B68").Value = "a1"
$ws.Range("B69").Value = "a2"
$ws.Range("B70").Value = "a3"
$ws.Range("B71").Value = "a4"
$ws.Range("B72").Value = "a5"
$ws.Range("B73").Value = "a6"
$ws.Range("B74").Value = "a7"
$ws.Range("B75").Value = "a8"
$ws.Range("B76").Value = "a9"
$ws.Range("B77").Value = "a888"
$ws.Range("B78").Value = "a1"
$ws.Range("B79").Value = "a1"
$ws.Range("B80").Value = "a2"
$ws.Range("B81").Value = "a888"
$ws.Range("B82").Value = "a1"
$ws.Range("B83").Value = "a0"

# The two cells that previously used a right-aligned numeric style now hold
# plain text, so drop the right alignment (reuse the ordinary wrap style).
$ws.Range("B14").HorizontalAlignment = 1
$ws.Range("B15").HorizontalAlignment = 1

# Column B no longer needs its custom width now that the values are shorter
# text tokens instead of numbers.
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668

# --- Switch the active sheet/selection from "survey" to "choices" ---
$ws.Activate()
$ws.Range("B6").Select()
